$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D look numeric (e.g. "1.00", "571.31") but must stay
# stored as text, exactly like the rest of the sheet (t="inlineStr"/shared
# string). Force a text number format on the whole price column before
# assigning values so Excel doesn't auto-convert them to numbers, then
# restore the default "Normal" style so no stray style index is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.539.86"
$ws.Range("E2").Value = "  -0.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.953.66"
$ws.Range("E3").Value = "  -1.92%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "571.31"
$ws.Range("E5").Value = "  -2.11%  "

# Row 6 - Solana
$ws.Range("D6").Value = "162.39"
$ws.Range("E6").Value = "  -0.28%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.16%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.949.80"
$ws.Range("E9").Value = "  -1.97%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  -3.20%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -4.40%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  +0.09%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  -2.87%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "34.86"
$ws.Range("E14").Value = "  +0.12%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.38%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "65.634.03"
$ws.Range("E16").Value = "  -0.44%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.443.84"

# Row 18 - Polkadot
$ws.Range("D18").Value = "7.09"
$ws.Range("E18").Value = "  +1.62%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "2.954.22"
$ws.Range("E19").Value = "  -1.96%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "15.95"
$ws.Range("E20").Value = "  +14.42%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "445.72"
$ws.Range("E21").Value = "  -2.61%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  +1.18%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -1.16%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "82.03"
$ws.Range("E24").Value = "  -0.64%  "

# Row 25 - Fetch.AI
$ws.Range("D25").Value = "2.25"
$ws.Range("E25").Value = "  -3.58%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "12.25"
$ws.Range("E26").Value = "  -1.06%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  -5.17%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  +0.05%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  +8.55%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "8.09"
$ws.Range("E30").Value = "  -0.56%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.52%  "

# Row 32 - PEPE
$ws.Range("E32").Value = "  -0.20%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +4.24%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "27.16"
$ws.Range("E34").Value = "  +0.28%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.07%  "

# Row 36 - Mantle
$ws.Range("E36").Value = "  -1.92%  "

# Row 37 - Filecoin
$ws.Range("E37").Value = "  -1.50%  "

# Row 38 - Arweave
$ws.Range("D38").Value = "46.29"
$ws.Range("E38").Value = "  +6.12%  "

# Row 39 - OKB
$ws.Range("D39").Value = "49.10"
$ws.Range("E39").Value = "  -1.41%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  -7.47%  "

# Row 41 - was Kaspa, now TheGraph
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.304"
$ws.Range("E41").Value = "  -1.92%  "

# Row 42 - was TheGraph, now Kaspa
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.122"
$ws.Range("E42").Value = "  +0.67%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -6.33%  "

# Row 44 - Cosmos
$ws.Range("D44").Value = "8.52"
$ws.Range("E44").Value = "  +0.44%  "

# Row 45 - Bittensor
$ws.Range("D45").Value = "384.52"
$ws.Range("E45").Value = "  -0.94%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -1.67%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.679.98"
$ws.Range("E47").Value = "  -4.12%  "

# Row 48 - Monero
$ws.Range("D48").Value = "132.77"
$ws.Range("E48").Value = "  -1.78%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 - InjectiveProtocol
$ws.Range("E50").Value = "  -0.42%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  +1.38%  "

# Restore default styling on the price column so no extra "@" number-format
# style lingers on any cell (matches original, unstyled price cells).
$ws.Range("D2:D51").Style = "Normal"
